$p = $ppt.ActivePresentation
$ds = $p.Designs
Write-Host "Count:" $ds.Count
$d1 = $ds.Item(1)
Write-Host "Index:" $d1.Index
Write-Host "Preserved:" $d1.Preserved
